$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(9).Delete()
$ws.Range("A8").Select()
